# Append a new daily-tracking row (row 49) to the bottom of the sheet,
# mirroring the existing rows: A=date text, B=weekday text, C/D=numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column A to be treated as text first so that a date-shaped string
# like "2025/10/02" is stored as a literal string (matching the other date
# cells in the column) instead of being auto-converted into a date serial
# number. Reset the cell style back to Normal afterwards so the new cell
# does not end up with a lingering text number-format applied to it (the
# existing data rows carry no explicit style either).
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "2025/10/02"
$ws.Range("A49").Style = "Normal"

$ws.Range("B49").Value = "木"
$ws.Range("C49").Value = 11
$ws.Range("D49").Value = 25
